# Updated date in report.
#
# This script reproduces, via the Word COM/OOXML object model, the same
# net edits that Word itself would have produced when a user:
#   1. Changed the report date from "October 19 2018" to "October 21 2018"
#      (which leaves Word's automatic "_GoBack" last-edit bookmark sitting
#      right after the newly typed "21").
#   2-4. A few nearby paragraphs end up with some of their same-format
#      adjacent runs coalesced together -- an incidental side effect of
#      Word's run-management when text near them is touched/re-saved.
#   5. The old "_GoBack" bookmark (previously at the end of the "variable
#      does not seem to reset." paragraph) is removed, since a document can
#      only have one "_GoBack" bookmark and it now lives at change 1.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "October 19 2018" -> "October 21 2018", splitting the run at
# the edited word and dropping Word's "_GoBack" bookmark right after it.
# ---------------------------------------------------------------------

$dateRng = $d.Content
$dateRng.Find.Execute("19", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null

# Pin the run boundary before "19" so the preceding " " run is left alone,
# and drop the "_GoBack" bookmark right after "19"/"21" -- exactly where
# Word leaves it after the user's last text edit.
$dayBookmarkRange = $d.Range($dateRng.Start, $dateRng.Start)
$d.Bookmarks.Add("ZZZDateEditStart", $dayBookmarkRange) | Out-Null

$yearBookmarkRange = $d.Range($dateRng.End, $dateRng.End)
$d.Bookmarks.Add("_GoBack", $yearBookmarkRange) | Out-Null

$dateRng.Text = "21"

$d.Bookmarks.Item("ZZZDateEditStart").Delete() | Out-Null

# ---------------------------------------------------------------------
# Change 2: " if " / "Rover pitch and Rover roll" / " are" coalesce into
# a single run, without disturbing the unrelated runs before/after it.
# ---------------------------------------------------------------------

function Coalesce-Range($rangeToTouch) {
    # Nudging one character at the end of the range (flip then restore)
    # is enough to make Word re-evaluate and merge same-format runs that
    # lie fully inside the range, while bookmarks pinned at the range's
    # start/end keep the merge from bleeding into neighboring runs.
    $startPos = $rangeToTouch.Start
    $endPos = $rangeToTouch.End

    $before = $d.Range($startPos, $startPos)
    $d.Bookmarks.Add("ZZZEditBoundaryBefore", $before) | Out-Null
    $after = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("ZZZEditBoundaryAfter", $after) | Out-Null

    $lastChar = $d.Range($endPos - 1, $endPos)
    $savedText = $lastChar.Text
    $lastChar.Text = "#"
    $lastChar2 = $d.Range($endPos - 1, $endPos)
    $lastChar2.Text = $savedText

    $d.Bookmarks.Item("ZZZEditBoundaryBefore").Delete() | Out-Null
    $d.Bookmarks.Item("ZZZEditBoundaryAfter").Delete() | Out-Null
}

$pitchRollRng = $d.Content
$pitchRollRng.Find.Execute(" if Rover pitch and Rover roll are", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0) | Out-Null
Coalesce-Range $pitchRollRng

# ---------------------------------------------------------------------
# Change 4: "The grading rubric is used to lay out the order of the
# report." and the trailing " " run coalesce into one run.
# (Handled before change 3 so range offsets found below stay valid.)
# ---------------------------------------------------------------------

$rubricRng = $d.Content
$rubricRng.Find.Execute("The grading rubric is used to lay out the order of the report. ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
Coalesce-Range $rubricRng

# ---------------------------------------------------------------------
# Change 3: the tab run, the long "A check is done..." run and the
# trailing " " run coalesce into one run.
# ---------------------------------------------------------------------

$visionRng = $d.Content
$visionRng.Find.Execute("A check is done on a narrow range of vision", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$visionStart = $visionRng.Start - 1   # include the tab just before the text
$trailRng = $d.Content
$trailRng.Find.Execute("needs more adjustment. ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$visionEnd = $trailRng.End + 1        # include the extra trailing space run

$visionFullRng = $d.Range($visionStart, $visionEnd)
Coalesce-Range $visionFullRng

Write-Output "Done."
